$wb = $excel.ActiveWorkbook

# ----------------------------------------------------------------------
# Sheet "VENTAS POR GRUPO": insert a new row for SOLIS OCAMPO DIMAS ABDON
# (splitting out of what used to be row 10) and insert a brand new
# salesperson MORA RODRIGUEZ BYRON RIQUELME before MOROCHO PLAZA SHIRLEY
# AURELIA. Net effect: a new row is inserted before the VEHINVER SA row
# (old row 11), and the old "MOROCHO PLAZA SHIRLEY AURELIA" (old row 9)
# name moves to being preceded by the new "MORA RODRIGUEZ..." row, while
# its data shifts to the row that used to belong to "SOLIS OCAMPO DIMAS
# ABDON" (old row 10), whose data/name in turn moves to the newly
# inserted row 11.
# ----------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")

# Row 9: rename MOROCHO PLAZA SHIRLEY AURELIA -> MORA RODRIGUEZ BYRON RIQUELME
$ws1.Range("B9").Value = "MORA RODRIGUEZ BYRON RIQUELME"

# Row 10: rename SOLIS OCAMPO DIMAS ABDON -> MOROCHO PLAZA SHIRLEY AURELIA,
# and its I column (LAVABOS) value moves down to the new row, so zero it here.
$ws1.Range("B10").Value = "MOROCHO PLAZA SHIRLEY AURELIA"
$ws1.Range("I10").Value = 0

# Insert a new row 11 (inherits formatting from row above) for
# SOLIS OCAMPO DIMAS ABDON, carrying the 43.86 LAVABOS value.
$ws1.Range("A11").EntireRow.Insert()
$ws1.Range("A11").Value = "OFICINA-CATAECSA"
$ws1.Range("B11").Value = "SOLIS OCAMPO DIMAS ABDON"
$ws1.Range("C11").Value = 0
$ws1.Range("D11").Value = 0
$ws1.Range("E11").Value = 0
$ws1.Range("F11").Value = 0
$ws1.Range("G11").Value = 0
$ws1.Range("H11").Value = 0
$ws1.Range("I11").Value = 43.86
$ws1.Range("J11").Value = 0
$ws1.Range("K11").Value = 0
$ws1.Range("L11").Value = 0
$ws1.Range("M11").Value = 0
$ws1.Range("N11").Value = 0
$ws1.Range("O11").Value = 0
$ws1.Range("P11").Value = 0
$ws1.Range("Q11").Value = 0
$ws1.Range("R11").Value = 0

# The totals row (old row 12, now row 13) text counters go from "de 10" to
# "de 11" since there are now 11 data rows.
$ws1.Range("C13").Value = "0 de 11"
$ws1.Range("D13").Value = "0 de 11"
$ws1.Range("E13").Value = "1 de 11"
$ws1.Range("F13").Value = "0 de 11"
$ws1.Range("G13").Value = "0 de 11"
$ws1.Range("H13").Value = "0 de 11"
$ws1.Range("I13").Value = "1 de 11"
$ws1.Range("J13").Value = "0 de 11"
$ws1.Range("K13").Value = "0 de 11"
$ws1.Range("L13").Value = "0 de 11"
$ws1.Range("M13").Value = "1 de 11"
$ws1.Range("N13").Value = "0 de 11"
$ws1.Range("O13").Value = "1 de 11"
$ws1.Range("P13").Value = "0 de 11"
$ws1.Range("Q13").Value = "0 de 11"
$ws1.Range("R13").Value = "0 de 11"

# ----------------------------------------------------------------------
# Sheet "VENTA MENSUAL": same restructuring as above, but with fewer
# columns (abril..julio, PRESUPUESTO).
# ----------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")

# Row 9: rename MOROCHO PLAZA SHIRLEY AURELIA -> MORA RODRIGUEZ BYRON RIQUELME
$ws2.Range("B9").Value = "MORA RODRIGUEZ BYRON RIQUELME"
$ws2.Range("E9").Value = 0

# Row 10: rename SOLIS OCAMPO DIMAS ABDON -> MOROCHO PLAZA SHIRLEY AURELIA
$ws2.Range("B10").Value = "MOROCHO PLAZA SHIRLEY AURELIA"
$ws2.Range("E10").Value = 159.03
$ws2.Range("F10").Value = 0

# Insert a new row 11 for SOLIS OCAMPO DIMAS ABDON carrying the 43.86 julio
# value.
$ws2.Range("A11").EntireRow.Insert()
$ws2.Range("A11").Value = "OFICINA-CATAECSA"
$ws2.Range("B11").Value = "SOLIS OCAMPO DIMAS ABDON"
$ws2.Range("C11").Value = 0
$ws2.Range("D11").Value = 0
$ws2.Range("E11").Value = 0
$ws2.Range("F11").Value = 43.86
$ws2.Range("G11").Value = 0

$wb.Save()
